$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds the last-changed date for every data row.
# Update it from serial date 45172 (2023-09-03) to 45175 (2023-09-06)
# for all data rows (rows 2 through 111).
$ws.Range("C2:C111").Value = 45175
